$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header values for columns F, G, H
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header style (bold, border, centered/top-aligned) from an existing header cell (E1)
$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill boolean FALSE values for rows 2-12 in columns F, G, H
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 6).Value = $false
    $ws.Cells.Item($r, 7).Value = $false
    $ws.Cells.Item($r, 8).Value = $false
}
